# Add the "2022-Q1" sheet (fund holdings detail) positioned right after "2021-Q4",
# and update the "总计" (summary) sheet with a new 2022-Q1 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet by duplicating the "2021-Q4" sheet so that
#    it inherits the same column layout / header styling / borders, then
#    overwrite its contents with the 2022-Q1 fund holdings data.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)

$q1 = $wb.Worksheets.Item(2)
$q1.Name = "2022-Q1"

# The copied sheet only has the index-column style (A2) for row 2; extend the
# same index-column formatting down to rows 3-5 for the additional fund rows.
$q1.Range("A2").Copy()
$q1.Range("A3:A5").PasteSpecial(-4122)

# Fund holdings data for 2022-Q1: index, code, name, scale, total position,
# position ratio, holding value (100M yuan), position rank.
$rows = @(
    @(0, "159758", "华夏中证红利质量ETF", "1.81", "99.16", "4.59", "0.0831", 3),
    @(1, "519677", "银河定投宝中证腾讯济安价值100A股指数", "2.74", "91.56", "1.22", "0.0334", 9),
    @(2, "009263", "华宝红利精选混合A", "0.46", "83.67", "1.14", "0.0052", 2),
    @(3, "010841", "华宝红利精选混合C", "0.16", "83.67", "1.14", "0.0018", 2)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $q1.Cells.Item($r, 1).Value = $data[0]

    # Columns that may look numeric (fund code, scale, ratios, values) must be
    # forced to text so they keep their original textual representation.
    $bCell = $q1.Cells.Item($r, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $data[1]

    $q1.Cells.Item($r, 3).Value = $data[2]

    $dCell = $q1.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $data[3]

    $eCell = $q1.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $data[4]

    $fCell = $q1.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $data[5]

    $gCell = $q1.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $data[6]

    $q1.Cells.Item($r, 8).Value = $data[7]
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q1 above the
#    existing 2021-Q4 row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Match the index-column style used on the row below.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.12

$total.Range("A3").Value = 1
